# Update weekly "Fruta / hortaliza" data: refresh rows 15-18 with new
# prices/varieties and append the displaced rows as new rows 19-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: Bing / Primera ---
$ws.Range("D15").Value = 44917
$ws.Range("K15").Value = "Bing"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = 5000
$ws.Range("O15").Value = 6000
$ws.Range("P15").Value = 5625
$ws.Range("Q15").Value = "$/bandeja 10 kilos"
$ws.Range("S15").Value = 562
$ws.Range("T15").Value = 10

# --- Row 16: Santina / Primera ---
$ws.Range("D16").Value = 44917
$ws.Range("K16").Value = "Santina"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 5000
$ws.Range("O16").Value = 6000
$ws.Range("P16").Value = 5500
$ws.Range("Q16").Value = "$/bandeja 10 kilos"
$ws.Range("S16").Value = 550
$ws.Range("T16").Value = 10

# --- Row 17: Rainier / Segunda ---
$ws.Range("D17").Value = 44908
$ws.Range("K17").Value = "Rainier"
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15600
$ws.Range("Q17").Value = "$/caja 10 kilos"
$ws.Range("S17").Value = 1560
$ws.Range("T17").Value = 10

# --- Row 18: Brooks / Primera ---
$ws.Range("D18").Value = 44532
$ws.Range("K18").Value = "Brooks"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 27000
$ws.Range("O18").Value = 28000
$ws.Range("P18").Value = 27500
$ws.Range("Q18").Value = "$/bandeja 12 kilos"
$ws.Range("S18").Value = 2292
$ws.Range("T18").Value = 12

# --- Row 19 (new): Sweet Heart / Segunda ---
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value = 44580
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103001
$ws.Range("J19").Value = "Cereza"
$ws.Range("K19").Value = "Sweet Heart"
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 7000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 7500
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 750
$ws.Range("T19").Value = 10

# --- Row 20 (new): Brooks / Segunda ---
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44571
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103001
$ws.Range("J20").Value = "Cereza"
$ws.Range("K20").Value = "Brooks"
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 8500
$ws.Range("O20").Value = 9000
$ws.Range("P20").Value = 8750
$ws.Range("Q20").Value = "$/bandeja 10 kilos"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 875
$ws.Range("T20").Value = 10

# Apply the same date style used by the other date cells in column D
$ws.Range("D19").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("D20").NumberFormat = $ws.Range("D18").NumberFormat
